$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cell = $t.Cell(1, 1)
$cell.Range.Find.Execute("72÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "35÷5=", 1) | Out-Null
$cell = $t.Cell(1, 2)
$cell.Range.Find.Execute("77÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "99÷7=", 1) | Out-Null
$cell = $t.Cell(1, 3)
$cell.Range.Find.Execute("51÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "35÷3=", 1) | Out-Null
$cell = $t.Cell(1, 4)
$cell.Range.Find.Execute("23÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "16÷3=", 1) | Out-Null
$cell = $t.Cell(1, 5)
$cell.Range.Find.Execute("46÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "98÷5=", 1) | Out-Null
$cell = $t.Cell(5, 1)
$cell.Range.Find.Execute("27÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "15÷3=", 1) | Out-Null
$cell = $t.Cell(5, 2)
$cell.Range.Find.Execute("77÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "65÷5=", 1) | Out-Null
$cell = $t.Cell(5, 3)
$cell.Range.Find.Execute("35÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "68÷9=", 1) | Out-Null
$cell = $t.Cell(5, 4)
$cell.Range.Find.Execute("26÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "94÷3=", 1) | Out-Null
$cell = $t.Cell(5, 5)
$cell.Range.Find.Execute("76÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "97÷3=", 1) | Out-Null
$cell = $t.Cell(9, 1)
$cell.Range.Find.Execute("74÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "12÷7=", 1) | Out-Null
$cell = $t.Cell(9, 2)
$cell.Range.Find.Execute("77÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "19÷2=", 1) | Out-Null
$cell = $t.Cell(9, 3)
$cell.Range.Find.Execute("16÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "48÷4=", 1) | Out-Null
$cell = $t.Cell(9, 4)
$cell.Range.Find.Execute("35÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "46÷7=", 1) | Out-Null
$cell = $t.Cell(9, 5)
$cell.Range.Find.Execute("76÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "79÷5=", 1) | Out-Null
$cell = $t.Cell(13, 1)
$cell.Range.Find.Execute("84÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "42÷5=", 1) | Out-Null
$cell = $t.Cell(13, 2)
$cell.Range.Find.Execute("60÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "45÷9=", 1) | Out-Null
$cell = $t.Cell(13, 3)
$cell.Range.Find.Execute("32÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "95÷7=", 1) | Out-Null
$cell = $t.Cell(13, 4)
$cell.Range.Find.Execute("47÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "29÷7=", 1) | Out-Null
$cell = $t.Cell(13, 5)
$cell.Range.Find.Execute("86÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "28÷7=", 1) | Out-Null
$cell = $t.Cell(17, 1)
$cell.Range.Find.Execute("51÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "11÷7=", 1) | Out-Null
$cell = $t.Cell(17, 2)
$cell.Range.Find.Execute("97÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "58÷9=", 1) | Out-Null
$cell = $t.Cell(17, 3)
$cell.Range.Find.Execute("81÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "31÷2=", 1) | Out-Null
$cell = $t.Cell(17, 4)
$cell.Range.Find.Execute("54÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "53÷6=", 1) | Out-Null
$cell = $t.Cell(17, 5)
$cell.Range.Find.Execute("57÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "40÷7=", 1) | Out-Null
